# Auto-applied scheduled market-data refresh for Leve profit sheets.
# For each affected row, H:N (currentAveragePrice, currentAveragePriceNQ,
# currentAveragePriceHQ, LevePriceNQ, LevePriceHQ, LeveProfitNQ, LeveProfitHQ)
# are refreshed from the latest market snapshot. Cells that no longer apply
# (profit switches sign / side) are cleared rather than left stale.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Range("H9").Value = 227.66667
$ws.Range("I9").Value = 178.57143
$ws.Range("J9").Value = 296.4
$ws.Range("K9").Value = 178.57143
$ws.Range("L9").Value = 296.4
$ws.Range("M9").Value = -9.571429999999992
$ws.Range("N9").Value = -634.4
# Row 138
$ws.Range("H138").Value = 2085593.5
$ws.Range("I138").Value = 894.575
$ws.Range("J138").Value = 3574664
$ws.Range("K138").Value = 2683.725
$ws.Range("L138").Value = 10723992
$ws.Range("M138").Value = 2456.275
$ws.Range("N138").Value = -10734272
# Row 141
$ws.Range("H141").Value = 1261.1082
$ws.Range("I141").Value = 701.6429000000001
$ws.Range("J141").Value = 3001.6667
$ws.Range("K141").Value = 2104.9287
$ws.Range("L141").Value = 9005.000100000001
$ws.Range("M141").Value = 3075.0713
$ws.Range("N141").Value = -19365.0001

$ws = $wb.Worksheets.Item("ARM")
# Row 23
$ws.Range("H23").Value = 18487.5
$ws.Range("J23").Value = 18487.5
$ws.Range("L23").Value = 18487.5
$ws.Range("N23").Value = -19005.5
# Row 25
$ws.Range("H25").Value = 5714.2856
$ws.Range("I25").Value = 1000
$ws.Range("J25").Value = 7600
$ws.Range("K25").Value = 1000
$ws.Range("L25").Value = 7600
$ws.Range("M25").Value = -598
$ws.Range("N25").Value = -8404
# Row 29
$ws.Range("H29").Value = 10000
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 10000
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 10000
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = -10616
# Row 30
$ws.Range("H30").Value = 7881.8
$ws.Range("I30").Value = 2509
$ws.Range("K30").Value = 2509
$ws.Range("M30").Value = -2359
# Row 35
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()
# Row 36
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()
# Row 37
$ws.Range("H37").Value = 29000
$ws.Range("J37").Value = 29000
$ws.Range("L37").Value = 29000
$ws.Range("N37").Value = -29546
# Row 101
$ws.Range("H101").Value = 45000
$ws.Range("J101").Value = 45000
$ws.Range("L101").Value = 45000
$ws.Range("N101").Value = -51490

$ws = $wb.Worksheets.Item("BSM")
# Row 23
$ws.Range("H23").Value = 35500
$ws.Range("J23").Value = 35500
$ws.Range("L23").Value = 35500
$ws.Range("N23").Value = -36066
# Row 24
$ws.Range("H24").Value = 8764.909
$ws.Range("I24").Value = 482.8
$ws.Range("K24").Value = 482.8
$ws.Range("M24").Value = -247.8
# Row 25
$ws.Range("H25").Value = 3055.7273
$ws.Range("I25").Value = 714.125
$ws.Range("J25").Value = 9300
$ws.Range("K25").Value = 714.125
$ws.Range("L25").Value = 9300
$ws.Range("M25").Value = -479.125
$ws.Range("N25").Value = -9770
# Row 29
$ws.Range("H29").Value = 10000
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 10000
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 10000
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = -10578
# Row 31
$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()
# Row 37
$ws.Range("H37").Value = 11746.571
$ws.Range("I37").Value = 3013
$ws.Range("J37").Value = 15240
$ws.Range("K37").Value = 3013
$ws.Range("L37").Value = 15240
$ws.Range("M37").Value = -2876
$ws.Range("N37").Value = -15514

$ws = $wb.Worksheets.Item("CRP")
# Row 29
$ws.Range("H29").Value = 1122
$ws.Range("I29").Value = 1122
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 1122
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -829
$ws.Range("N29").ClearContents()
# Row 35
$ws.Range("H35").Value = 20000
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 20000
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 20000
$ws.Range("M35").ClearContents()
$ws.Range("N35").Value = -20588
# Row 36
$ws.Range("H36").Value = 15000
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 15000
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 15000
$ws.Range("M36").ClearContents()
$ws.Range("N36").Value = -15776
# Row 40
$ws.Range("H40").Value = 15000
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 15000
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 15000
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -15320
# Row 86
$ws.Range("H86").Value = 3593158.8
$ws.Range("I86").Value = 5057907
$ws.Range("J86").Value = 12662.333
$ws.Range("K86").Value = 5057907
$ws.Range("L86").Value = 12662.333
$ws.Range("M86").Value = -5056784
$ws.Range("N86").Value = -14908.333
# Row 89
$ws.Range("H89").Value = 3593158.8
$ws.Range("I89").Value = 5057907
$ws.Range("J89").Value = 12662.333
$ws.Range("K89").Value = 25289535
$ws.Range("L89").Value = 63311.665
$ws.Range("M89").Value = -25283919
$ws.Range("N89").Value = -74543.66500000001

$ws = $wb.Worksheets.Item("CUL")
# Row 23
$ws.Range("H23").Value = 112
$ws.Range("I23").Value = 63.2
$ws.Range("J23").Value = 139.11111
$ws.Range("K23").Value = 189.6
$ws.Range("L23").Value = 417.33333
$ws.Range("M23").Value = 45.39999999999998
$ws.Range("N23").Value = -887.3333299999999
# Row 24
$ws.Range("H24").Value = 209.2
$ws.Range("I24").Value = 209.2
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 627.5999999999999
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = -397.5999999999999
$ws.Range("N24").ClearContents()
# Row 25
$ws.Range("H25").Value = 1749.6364
$ws.Range("I25").Value = 149.5
$ws.Range("J25").Value = 2105.2222
$ws.Range("K25").Value = 448.5
$ws.Range("L25").Value = 6315.6666
$ws.Range("M25").Value = -279.5
$ws.Range("N25").Value = -6653.6666
# Row 30
$ws.Range("H30").Value = 1749.6364
$ws.Range("I30").Value = 149.5
$ws.Range("J30").Value = 2105.2222
$ws.Range("K30").Value = 448.5
$ws.Range("L30").Value = 6315.6666
$ws.Range("M30").Value = -346.5
$ws.Range("N30").Value = -6519.6666
# Row 131
$ws.Range("H131").Value = 920.87
$ws.Range("J131").Value = 921.7879
$ws.Range("L131").Value = 2765.3637
$ws.Range("N131").Value = -12845.3637

$ws = $wb.Worksheets.Item("GSM")
# Row 51
$ws.Range("H51").Value = 29900
$ws.Range("J51").Value = 29900
$ws.Range("L51").Value = 29900
$ws.Range("N51").Value = -30918
# Row 80
$ws.Range("H80").Value = 2182.8333
$ws.Range("I80").Value = 2185.7144
$ws.Range("J80").Value = 2178.8
$ws.Range("K80").Value = 2185.7144
$ws.Range("L80").Value = 2178.8
$ws.Range("M80").Value = -1187.7144
$ws.Range("N80").Value = -4174.8
# Row 83
$ws.Range("H83").Value = 2182.8333
$ws.Range("I83").Value = 2185.7144
$ws.Range("J83").Value = 2178.8
$ws.Range("K83").Value = 10928.572
$ws.Range("L83").Value = 10894
$ws.Range("M83").Value = -5936.572
$ws.Range("N83").Value = -20878
# Row 113
$ws.Range("H113").Value = 1464.5
$ws.Range("I113").Value = 1298.75
$ws.Range("J113").Value = 1685.5
$ws.Range("K113").Value = 1298.75
$ws.Range("L113").Value = 1685.5
$ws.Range("M113").Value = 871.25
$ws.Range("N113").Value = -6025.5

$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 799.8
$ws.Range("I96").Value = 666.6667
$ws.Range("J96").Value = 999.5
$ws.Range("K96").Value = 666.6667
$ws.Range("L96").Value = 999.5
$ws.Range("M96").Value = 706.3333
$ws.Range("N96").Value = -3745.5
# Row 126
$ws.Range("H126").Value = 892.7917
$ws.Range("I126").Value = 818.7059
$ws.Range("K126").Value = 2456.1177
$ws.Range("M126").Value = 13.88229999999976
# Row 132
$ws.Range("H132").Value = 2832.5715
$ws.Range("I132").Value = 3526.32
$ws.Range("J132").Value = 1812.3529
$ws.Range("K132").Value = 10578.96
$ws.Range("L132").Value = 5437.0587
$ws.Range("M132").Value = -8048.960000000001
$ws.Range("N132").Value = -10497.0587
# Row 136
$ws.Range("H136").Value = 3296.0833
$ws.Range("I136").Value = 4222.1875
$ws.Range("J136").Value = 2555.2
$ws.Range("K136").Value = 12666.5625
$ws.Range("L136").Value = 7665.599999999999
$ws.Range("M136").Value = -10116.5625
$ws.Range("N136").Value = -12765.6
